# R22 UAT2 - Regression: update CDRCancellation reference number and add a
# new "New Data" sheet ahead of it, becoming the active tab.

$wb = $excel.ActiveWorkbook

# 1) Update the existing CDRCancellation sheet's reference value first so
#    the new shared-string entry lands at the lower index.
$cdrSheet = $wb.Worksheets.Item("CDRCancellation")
$cdrSheet.Range("A2").Value = "CLK0601333"

# 2) Insert a brand-new sheet ahead of the active (CDRCancellation) sheet,
#    rename it, and populate it with its own reference row.
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "New Data"
$newSheet.Range("A1").Value = "CREDIT.THEIR.REF"
$newSheet.Range("A2").Value = "CLK0601322"
